# Generate Report for Handoff
# Updates the "Ready for handoff" rows (rows 4-7) on both the zh-cn and
# de-de localization-status sheets: Priority moves from "low" to "ht"
# (handed-off) and the Latest Handoff Datetime is refreshed to reflect
# the new handoff run.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zhHandoffTime = "2016-08-19 02:30:16"
foreach ($r in 4..7) {
    $zh.Cells.Item($r, 5).Value = "ht"
    $zh.Cells.Item($r, 8).Value = $zhHandoffTime
}

$de = $wb.Worksheets.Item("de-de")
$deHandoffTime = "2016-08-19 02:30:21"
foreach ($r in 4..7) {
    $de.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 8).Value = $deHandoffTime
}

# The Overview sheet's "Latest HO Xliff Generate Date" column shares the
# same string value as the de-de handoff datetime, so it picks up the
# refreshed timestamp automatically once the shared string is updated.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = $deHandoffTime
}
